# Apply the settings.xlsx updates described in the commit:
#  - switch the dev project panel name to the smaller panel
#  - drop the HC-4 anchor id, keep only HC-100
#  - turn normalization and database injection steps on
#  - leave the active selection on B8 (last cell touched by the user)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# project_name: dev_project_panel_1 -> small_dev_project_panel_1
$ws.Range("B3").Value = "small_dev_project_panel_1"

# do_normalization: 0 -> 1
$ws.Range("B7").Value = 1

# anchor_ids: "HC-4, HC-100" -> "HC-100"
$ws.Range("B8").Value = "HC-100"

# do_database_injection: 0 -> 1
$ws.Range("B10").Value = 1

# Reflect the final user selection on the sheet
$ws.Range("B8").Select()

$wb.Save()
